$wb = $excel.ActiveWorkbook

# Both the "展览" sheet (1) and the "全部类型" sheet (4) contain identical data
# and receive identical edits in this revision.
$sheetIndexes = @(1, 4)

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)

    # --- Refresh the "想去人数" / "最低票价" counters that simply incremented. ---
    $ws.Range("F2").Value = 270
    $ws.Range("F3").Value = 1385
    $ws.Range("F6").Value = 236
    $ws.Range("F8").Value = 16
    $ws.Range("F9").Value = 187
    $ws.Range("F10").Value = 137
    $ws.Range("F11").Value = 4658
    $ws.Range("F12").Value = 6926
    $ws.Range("F16").Value = 574
    $ws.Range("F18").Value = 4153
    $ws.Range("F19").Value = 802
    $ws.Range("F20").Value = 77
    $ws.Range("F21").Value = 64
    $ws.Range("G21").Value = 50
    $ws.Range("F22").Value = 2738
    $ws.Range("F23").Value = 574
    $ws.Range("F25").Value = 171
    $ws.Range("F27").Value = 376
    $ws.Range("F29").Value = 236
    $ws.Range("F31").Value = 1640
    $ws.Range("F32").Value = 1038
    $ws.Range("F33").Value = 69
    $ws.Range("F34").Value = 428

    # --- Insert the new "赣州·十万伏特" row above row 37, pushing the rest down. ---
    $ws.Rows.Item(37).Insert()

    # Row insertion can copy neighbouring formatting onto the new row; make sure
    # the index cell A37 keeps the same bold/bordered/centered look as every
    # other row in column A.
    $a37 = $ws.Range("A37")
    $a37.Font.Bold = $true
    $a37.HorizontalAlignment = -4108
    $a37.VerticalAlignment = -4160
    $a37.Borders.LineStyle = 1

    $ws.Range("A37").Value = 36
    # Leading apostrophe forces text, matching the existing date-as-text cells
    # in column B (otherwise Excel reinterprets the string as a date serial).
    # Resetting the style afterwards drops the quote-prefix formatting so the
    # cell ends up with the same default (unstyled) look as its neighbours.
    $ws.Range("B37").Value = "'2024-08-10"
    $ws.Range("B37").Style = "Normal"
    $ws.Range("C37").Value = "赣州·十万伏特-星铁&音乐 次元音乐only2.0"
    $ws.Range("D37").Value = "平安大道 麋鹿LiveHouse"
    $ws.Range("E37").Value = "2024.08.10 10:00-08.10 17:00"
    $ws.Range("F37").Value = 0
    $ws.Range("G37").Value = 45
    $ws.Range("H37").Value = "https://show.bilibili.com/platform/detail.html?id=89411"
    $ws.Range("I37").Value = "//i1.hdslb.com/bfs/openplatform/202407/H7fGQbYD1721132195795.jpeg"

    # The "南昌·CM03" listing (now on row 41 after the shift) also had its
    # counter bumped from 171 to 181 in this same refresh.
    $ws.Range("F41").Value = 181
}
